$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New TPM-derived values for the Il33-Il1rl1 sending/target cluster pairs.
# The "ECs" cluster is dropped as a possible *target* (column D), and the
# "Resolving-Mac" cluster is dropped as a possible *source* row-block, which
# shrinks the data from 16 rows (4x4 clusters) to 12 rows (4x3).
$rows = @(
        @("ECs","Il33","Il1rl1","FAPs",1,0.3333333333333333,0.1829693333333333,0.548908,0.004221684928116883,0.004221684928116882,3,1,8.061744666666668,24.185234,0.8690987204662018,0.8690987204662017,1.475052047163556,13.275468424472,0.003669060969237832,0.003669060969237831),
        @("ECs","Il33","Il1rl1","MuSCs",1,0.3333333333333333,0.1829693333333333,0.548908,0.004221684928116883,0.004221684928116882,1,0.3333333333333333,0.086217,0.258651,0.009294648674778319,0.009294648674778319,0.015775067012,0.141975603108,[double]"3.923907822245319E-05",[double]"3.923907822245318E-05"),
        @("ECs","Il33","Il1rl1","Resolving-Mac",1,0.3333333333333333,0.1829693333333333,0.548908,0.004221684928116883,0.004221684928116882,3,1,1.128021,3.384063,0.1216066308590198,0.1216066308590198,0.206393250356,1.857539253204,0.0005133848806565974,0.0005133848806565973),
        @("FAPs","Il33","Il1rl1","FAPs",3,1,42.68588366666666,128.057651,0.9848992092604805,0.9848992092604804,3,1,8.061744666666668,24.185234,0.8690987204662018,0.8690987204662017,344.1226949917038,3097.104254925334,0.8559746425564575,0.8559746425564574),
        @("FAPs","Il33","Il1rl1","MuSCs",3,1,42.68588366666666,128.057651,0.9848992092604805,0.9848992092604804,1,0.3333333333333333,0.086217,0.258651,0.009294648674778319,0.009294648674778319,3.680248832089,33.122239488801,0.009154292130143139,0.009154292130143139),
        @("FAPs","Il33","Il1rl1","Resolving-Mac",3,1,42.68588366666666,128.057651,0.9848992092604805,0.9848992092604804,3,1,1.128021,3.384063,0.1216066308590198,0.1216066308590198,48.15057317955699,433.3551586160129,0.1197702745738798,0.1197702745738798),
        @("MuSCs","Il33","Il1rl1","FAPs",1,0.3333333333333333,0.1203946666666667,0.361184,0.002777888187231682,0.002777888187231682,3,1,8.061744666666668,24.185234,0.8690987204662018,0.8690987204662017,0.9705910618951112,8.735319557056,0.002414259069121232,0.002414259069121232),
        @("MuSCs","Il33","Il1rl1","MuSCs",1,0.3333333333333333,0.1203946666666667,0.361184,0.002777888187231682,0.002777888187231682,1,0.3333333333333333,0.086217,0.258651,0.009294648674778319,0.009294648674778319,0.010380066976,0.09342060278400001,[double]"2.58194947581353E-05",[double]"2.58194947581353E-05"),
        @("MuSCs","Il33","Il1rl1","Resolving-Mac",1,0.3333333333333333,0.1203946666666667,0.361184,0.002777888187231682,0.002777888187231682,3,1,1.128021,3.384063,0.1216066308590198,0.1216066308590198,0.135807712288,1.222269410592,0.0003378096233523149,0.0003378096233523149),
        @("Resolving-Mac","Il33","Il1rl1","FAPs",3,1,0.3511096666666667,1.053329,0.00810121762417095,0.008101217624170948,3,1,8.061744666666668,24.185234,0.8690987204662018,0.8690987204662017,2.830556482665112,25.475008343986,0.007040757871385217,0.007040757871385214),
        @("Resolving-Mac","Il33","Il1rl1","MuSCs",3,1,0.3511096666666667,1.053329,0.00810121762417095,0.008101217624170948,1,0.3333333333333333,0.086217,0.258651,0.009294648674778319,0.009294648674778319,0.030271622131,0.272444599179,[double]"7.529797165459128E-05",[double]"7.529797165459127E-05"),
        @("Resolving-Mac","Il33","Il1rl1","Resolving-Mac",3,1,0.3511096666666667,1.053329,0.00810121762417095,0.008101217624170948,3,1,1.128021,3.384063,0.1216066308590198,0.1216066308590198,0.3960590773029999,3.564531695727,0.0009851617811311424,0.0009851617811311422)
)

$nRows = $rows.Count
$nCols = $rows[0].Count
$data = New-Object 'object[,]' $nRows, $nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $nRows - 1
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $nCols)).Value = $data

# The old sheet had 16 data rows (rows 2-17); only 12 remain now, so drop the
# trailing rows entirely (shrinks dimension from A1:T17 to A1:T13).
$lastOldRow = 17
if ($lastOldRow -gt $endRow) {
    $ws.Range("A$($endRow + 1):T$lastOldRow").EntireRow.Delete() | Out-Null
}
